$d = $word.ActiveDocument
$wdReplaceOne = 1
$results = @()

$old1 = "In this project, we will be adding a fast search feature to the resources section of this application. We are interpreting this as picking a type of resource, choosing an attribute to search by, and typing in what we want to search. Then, in the calendar view, only events containing that resource-attribute match will be displayed. This is very similar to the filter feature, except you have to uncheck every single resource you don’t want to see if you only want to view one. Our addition will make displaying a single resource much easier."
$new1 = "In this project, we will be adding a fast search feature to the resources section of Rapla. We are interpreting this as: picking a type of resource, choosing an attribute to search by, and typing in search criterion. Then, in the calendar view, only events containing that resource-attribute match will be displayed. This is very similar to the filter feature, except filter requires unchecking every single resource you don’t want to see (if you only want to view one kind of resource). Our addition will make filtering/searching for a single resource much easier. We will refer to it as SearchEditButton."
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
$results += $found1

$old2 = "Along with this, we would like to use the class ClassifibleFilterEdit located in org.rapla.gui.internal.edit. This class seems to help create text fields and combo boxes inside the popup a filter button would have, which helps in our search feature. We would like our user to choose a resource to search through a combo box, "
$new2 = "Along with this, we would like to use the class ClassifibleFilterEdit located in org.rapla.gui.internal.edit. This class seems to help create TextFields and JComboBoxes inside the popup a filter button would have, which would also be needed in our SearchEditButton. We would like our user to choose a resource to search by using a JComboBox, "
$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
$results += $found2

$old3 = "then an attribute through a combo box, then type in their criteria through a text field."
$new3 = "then an attribute through another JComboBox, then type in their criteria through a TextField."
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
$results += $found3

$old4 = "We are also thinking about using ClassificationFilterRuleImpl & ClassificationFilterImpl from org.rapla.dynamictype.internal or ClassificationFilterRule & ClassificationFilter from org.rapla.dynamictype. These parts of the program seem to contain the code that uses filter rules that choose which events in the Rapla interface fit defined rules. In this case, we can manipulate it to match events with search criteria."
$new4 = "We are also thinking about using ClassificationFilterRuleImpl & ClassificationFilterImpl from org.rapla.dynamictype.internal. These parts of the program seem to contain the code that uses filter rules that choose which events in the Rapla interface are displayed. In this case, we can manipulate it to work with SearchEditButton also."
$found4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
$results += $found4

$old5 = "Lastly, we will need to change the code a bit in the section of this project that actually adds each part of the program to the interface frame. This is simply where we will add the search arrow button."
$new5 = "Lastly, we will need to change the code a bit in the section of this project that actually adds buttons to the interface frame. This is simply where we will add SearchEditButton to a certain spot in the GUI."
$found5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
$results += $found5

$old6 = "The existing code helped a lot in our design process. Our first gut reaction to working on this assignment at first was to explore the program and see what features already existed, from a user perspective. Once we found how closely resembling the filter feature was to search, we knew we should change that around to make our work progress well. This way, we won’t have to develop any algorithms for searching and just use ones that already work perfectly fit with the program."
$new6 = "The existing code helped a lot in our design process. Our first gut reaction to working on this assignment was to explore Rapla and see what features already existed, from a user perspective. Once we found how closely resembling the filter feature was to a search feature, we knew we should specialize it to make our addition. This way, we won’t have to develop any new algorithms for searching and just use ones that already fit the program."
$found6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
$results += $found6

$old7 = "Also, the project was designed in such a way that most of the code is generalized and can be reused in many different cases. This makes it very easy for extra add-ons to be implemented, as it seems designed to be easy for others to collaborate on."
$new7 = "Also, the project was designed in such a way that most of the code is generalized and can be reused in many different cases. This makes it very easy for extra add-ons to be implemented, as it seems designed to be easy for others to collaborate on. For example, RaplaArrowButton simply allows the user to create a box with any text, but the box contains an arrow that suggests a drop-down box to appear when clicked. This class can be extended to others to create many different types of RaplaArrowButtons."
$found7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
$results += $found7

$old8 = "At this point in our design, we would like to design a class that uses bits from each of the classes mentioned in 1. "
$new8 = "At this point in our design, we would like to design our UML to map out in a formal way how SearchEditButton will relate to all the other classes previously mentioned. "
$found8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
$results += $found8

$old9 = "the parent class of a diagram is located at the top. "
$new9 = "the parent at the highest level of hierarchy is located at the top of the diagram, and its subsequent children are immediately below (and so on and so forth). "
$found9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
$results += $found9

$old10 = "Since our search button will branch from the filter button, the filter button will be located above the search button."
$new10 = "Since SearchEditButton will branch from the FilterEditButton, FilterEditButton will be located above SearchEditButton."
$found10 = $d.Content.Find.Execute($old10, $true, $false, $false, $false, $false, $true, 1, $false, $new10, 2)
$results += $found10

$old11 = "Then, ClassifiableFilterEdit and RaplaArrowButton are in between the filter button and search button. That is because both buttons aggregate with these classes, and this placement prevents any arrows from hitting each other."
$new11 = "Then, ClassifiableFilterEdit and RaplaArrowButton are in between FilterEditButton and SearchEditButton. That is because both buttons aggregate with these classes, and this placement prevents any arrows from intersecting."
$found11 = $d.Content.Find.Execute($old11, $true, $false, $false, $false, $false, $true, 1, $false, $new11, 2)
$results += $found11

$old12 = "Then there is the search button and the two Impl classes that we would like to use as a part of our code. Since we want to call on the filter functions included in these Impl classes, we use a dependency arrow to show that our search button depends on these two classes. Since this relation is independent from the filter button, we put the Impl classes at the very bottom to show that their significance only exists with our search button."
$new12 = "Then there are ClassificationFilterRuleImpl and ClassificationFilterImpl. We want to call on functions included in these classes, so we use a dependency arrow to show that SearchEditButton depends on these two classes. Since this relation is independent from FilterEditButton, we put these classes at the very bottom to show that their significance only exists with SearchEditButton."
$found12 = $d.Content.Find.Execute($old12, $true, $false, $false, $false, $false, $true, 1, $false, $new12, 2)
$results += $found12

$old13 = "First, our diagram shows that SearchEditButton extends FilterEditButton. We want to inherit all methods within the FilterEditButton because SearchEditButton will work the same way; has an arrow, we click on it, a popup shows, and you choose/type search criteria. Some specialization will have to take place, so we will probably end up editing the constructor method a bit in SearchEditButton."
$new13 = "First, our diagram shows that SearchEditButton extends FilterEditButton. We want to inherit all methods and private variables within FilterEditButton because SearchEditButton will work the same way; has an arrow, we click on it, a popup shows, and you choose/type search criteria. Some specialization will have to take place, so we will probably end up editing the constructor method a bit in SearchEditButton."
$found13 = $d.Content.Find.Execute($old13, $true, $false, $false, $false, $false, $true, 1, $false, $new13, 2)
$results += $found13

$old14 = "Next, we use aggregation to show that SearchEditButton has ClassifiableFilterEdit and RaplaArrowButton. It contains an instance of RaplaArrowButton, because the button has a literal arrow printed on it and we will use this button to open/close the search popup. It has an instance of ClassifiableFilterEdit because this class contains TextFields and JComboBoxes; these will be used to enter user input and choose search attributes (respectively)."
$new14 = "Next, we use aggregation to show that SearchEditButton has RaplaArrowButton and ClassifiableFilterEdit. It contains an instance of RaplaArrowButton, because the button has a literal arrow printed on it and we will use this button to open/close the search popup. It has an instance of ClassifiableFilterEdit because this class contains TextFields and JComboBoxes; these will be used to enter user input and choose search attributes (respectively)."
$found14 = $d.Content.Find.Execute($old14, $true, $false, $false, $false, $false, $true, 1, $false, $new14, 2)
$results += $found14

$old15 = "related to filtering resource types; it checks events in the current view for a filtered resource type and only shows events that meet this criterion. Since search works basically the same way, we can use this algorithm process to implement our feature. We only need to call on these methods, we use dependency."
$new15 = "related to filtering resource types; it checks events in the current view for a filtered resource type and only shows events that meet this criterion. Since search works basically the same way, we can use this algorithm process to implement SearchEditButton. We only need to call on these methods, so we use dependency."
$found15 = $d.Content.Find.Execute($old15, $true, $false, $false, $false, $false, $true, 1, $false, $new15, 2)
$results += $found15

$old16 = "Here’s our UML diagram: The .cld and .jpg file are both separately included in our repository, because this is probably difficult to view."
$new16 = "Here’s our UML diagram: The .cld and .jpg file are both separately included in our repository, because this is probably difficult to view in this file."
$found16 = $d.Content.Find.Execute($old16, $true, $false, $false, $false, $false, $true, 1, $false, $new16, 2)
$results += $found16

$old17 = "And here is the screenshot of us committing our files. See the message to find their location."
$new17 = "And here is the screenshot of us committing our files. We will organize all our files to be located in a folder titles “assignment related documents” so all our specific contributions stay independent from files the Rapla developers created."
$found17 = $d.Content.Find.Execute($old17, $true, $false, $false, $false, $false, $true, 1, $false, $new17, 2)
$results += $found17

Write-Output ($results -join ",")
